$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn on AutoFilter for the data range (adds the hidden _FilterDatabase
# defined name scoped to this sheet, matching Data > Filter in the UI).
$flt = $ws.Names.Add("_xlnm._FilterDatabase", "='Top 50'!`$A`$1:`$F`$51")
$flt.Visible = $false

# Relabel the "Sport" column (column C) values:
#   Football -> American Football   (do this first so we don't double-rename)
#   Soccer   -> Football
for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "Football") {
        $cell.Value = "American Football"
    }
}
for ($r = 2; $r -le 51; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq "Soccer") {
        $cell.Value = "Football"
    }
}

# Move the active selection to C3 (matches the saved view in the workbook).
$ws.Range("C3").Select()
